$wb = $excel.ActiveWorkbook

# The existing (first) worksheet in the workbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1 and name it "Sheet2"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate values A1:A5
$ws2.Range("A1").Value = 1
$ws2.Range("A2").Value = 2
$ws2.Range("A3").Value = 3
$ws2.Range("A4").Value = 4
$ws2.Range("A5").Value = 5

# Formulas in A7:A11 (row 6 intentionally left blank)
$ws2.Range("A7").Formula = "=SUM(A1:A5)"
$ws2.Range("A8").Formula = "=AVERAGE(A1:A5)"
$ws2.Range("A9").Formula = "=MAX(A1:A5)"
$ws2.Range("A10").Formula = "=MIN(A1:A5)"
$ws2.Range("A11").Formula = "=SUBTOTAL(6,A1:A5)"

# Sheet2 becomes the active sheet/tab, with the entire row 11 selected
$ws2.Activate()
$ws2.Range("A11:XFD11").Select()

$wb.Save()
